$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New practice records for contests ABC139 (questions A-E) and ABC129 (questions D, E)
$date1 = Get-Date -Year 2019 -Month 9 -Day 1
$date2 = Get-Date -Year 2019 -Month 9 -Day 2

# Row 122: ABC139 A
$ws.Cells.Item(122, 1).Value = 139
$ws.Cells.Item(122, 2).Value = "A"
$ws.Cells.Item(122, 3).Value = "AC"
$ws.Cells.Item(122, 5).Value = $date1
$ws.Cells.Item(122, 5).NumberFormat = "m/d/yyyy"

# Row 123: ABC139 B
$ws.Cells.Item(123, 1).Value = 139
$ws.Cells.Item(123, 2).Value = "B"
$ws.Cells.Item(123, 3).Value = "AC"
$ws.Cells.Item(123, 5).Value = $date1
$ws.Cells.Item(123, 5).NumberFormat = "m/d/yyyy"

# Row 124: ABC139 C
$ws.Cells.Item(124, 1).Value = 139
$ws.Cells.Item(124, 2).Value = "C"
$ws.Cells.Item(124, 3).Value = "AC"
$ws.Cells.Item(124, 5).Value = $date1
$ws.Cells.Item(124, 5).NumberFormat = "m/d/yyyy"

# Row 125: ABC139 D
$ws.Cells.Item(125, 1).Value = 139
$ws.Cells.Item(125, 2).Value = "D"
$ws.Cells.Item(125, 3).Value = "AC"
$ws.Cells.Item(125, 5).Value = $date1
$ws.Cells.Item(125, 5).NumberFormat = "m/d/yyyy"

# Row 126: ABC139 E
$ws.Cells.Item(126, 1).Value = 139
$ws.Cells.Item(126, 2).Value = "E"
$ws.Cells.Item(126, 3).Value = "AC"
$ws.Cells.Item(126, 4).Value = $true
$ws.Cells.Item(126, 5).Value = $date2
$ws.Cells.Item(126, 5).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(126, 6).Value = "dfs, dag"

# Row 127: ABC129 D
$ws.Cells.Item(127, 1).Value = 129
$ws.Cells.Item(127, 2).Value = "D"
$ws.Cells.Item(127, 3).Value = "AC"
$ws.Cells.Item(127, 5).Value = $date2
$ws.Cells.Item(127, 5).NumberFormat = "m/d/yyyy"

# Row 128: ABC129 E
$ws.Cells.Item(128, 1).Value = 129
$ws.Cells.Item(128, 2).Value = "E"
$ws.Cells.Item(128, 3).Value = "AC"
$ws.Cells.Item(128, 5).Value = $date2
$ws.Cells.Item(128, 5).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(128, 6).Value = "digit dp"

# Update frozen pane / view to match the new extent
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.FreezePanes = $true
$ws.Range("A120").Select()
$excel.ActiveWindow.ScrollRow = 120
$ws.Range("G129").Select()
